$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LuSTRE")

# Row 3 becomes the bif:contains (virtuoso) variant of the query that used
# to live in row 5; rows 4 and 5 (the skos:definition-limited and
# rdfs:label-limited variants) are removed entirely.

$querystring = @"
SELECT DISTINCT ?subject ?predicate ?object
WHERE {
  ?subject rdfs:label ?object .
  ?subject ?predicate ?object .
  ?object bif:contains "!!!SEARCHWORD!!!" 
  OPTION (score ?sc) .
} 
ORDER BY DESC (?sc)
LIMIT 100
"@

$description = @"
Return subjects whose rdfs:label bif:contains searchword; ordered by score and limited to 100.
# http://www.openlinksw.com/schemas/bif# is a feature of SPARQL Virtuoso server, see http://docs.openlinksw.com/virtuoso/rdfsparqlrulefulltext/
# could be made more permissive with wildcard " ' !!!SEARCHWORD!!! * ' " but also much slower
"@

$ws.Range("C3").Value = "virtuoso"
$ws.Range("D3").Value = $querystring
$ws.Range("E3").Value = $description
$ws.Range("F3").Value = 0

# Remove the now-redundant rows 4 and 5 (shifts everything below up).
$ws.Rows("4:5").Delete()

# Row 3 grows taller to fit the longer bif:contains description/query text.
$ws.Rows("3").RowHeight = 135

$ws.Range("A4:XFD24").Select()
